# Automatische test-sync: 2025-08-05 18:43:50
#
# Adds a new test-mail row (row 39) to the "Logs" sheet, extends the
# conditional formatting ranges to cover it, and updates the "Dashboard"
# summary sheet so the "Inkoop / Bestellingen" category (now 5) sorts
# above "Klantenservice / Contact" (still 4).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append the new row of data
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A39").Value = "Bestel je 200 stuks M8-bouten RVS voor Van Dijk?"
$logs.Range("B39").Value = "mailmind.test@zohomail.eu"
$logs.Range("C39").Value = "Testmail #18: Bestel je 200 stuks M8-bouten RVS voor Van Dijk?"
$logs.Range("D39").Value = "Inkoop / Bestellingen"
$logs.Range("E39").Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@bedrijf.nl."
$logs.Range("F39").Value = "2025-08-05 18:43:46"
$logs.Range("G39").Value = "Ja"
$logs.Range("H39").Value = "Ja"
$logs.Range("I39").Value = "Nee"
$logs.Range("J39").Value = "Nee"

# ---------------------------------------------------------------------
# 2. Logs sheet: extend the conditional-formatting ranges from row 38
#    down to row 39 (D, G, H, I, J columns)
# ---------------------------------------------------------------------
$ccols = "D", "G", "H", "I", "J"
foreach ($col in $ccols) {
    $oldRange = $logs.Range($col + "2:" + $col + "38")
    $newRange = $logs.Range($col + "2:" + $col + "39")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 3. Dashboard sheet: swap rows 3 and 4 so the category counts stay
#    sorted descending now that "Inkoop / Bestellingen" grew to 5
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "Inkoop / Bestellingen"
$dash.Range("B3").Value = 5
$dash.Range("A4").Value = "Klantenservice / Contact"
$dash.Range("B4").Value = 4
